# Append one new row (row 93) of portfolio data for 2025-11-16,
# mirroring the values already present on the prior row (92), per the
# commit "Update portfolio-updates.xlsx on 2025-11-16 12:57:32".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A stores the date as plain text (e.g. "2025-11-15" on row 92),
# not as a real Excel date. Force the cell to Text format first so the
# COM layer doesn't auto-convert the "YYYY-MM-DD" literal into a date
# serial number, then restore the default (Normal) style so the cell
# doesn't end up carrying a leftover number format.
$ws.Range("A93").NumberFormat = "@"
$ws.Range("A93").Value = "2025-11-16"
$ws.Range("A93").Style = "Normal"

# Columns B/C/D are plain numeric values, same as row 92.
$ws.Range("B93").Value = 57.68000030517578
$ws.Range("C93").Value = 391.2000122070312
$ws.Range("D93").Value = 303.75
